$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new "riderID" column at B (shifts old B:D phoneNumber/totalBagVolume/currentAvailableBagVolume to C:E)
$ws.Columns.Item(2).Insert()

# --- Rider names in column A (rows 2-9).
# Entered in this specific order so the shared-string table is built in the
# same sequence as the authored workbook (Harsh is typed before Jagdish,
# then rows 6/7 are corrected to their final order).
$ws.Cells.Item(2, 1).Value = "Shree"
$ws.Cells.Item(3, 1).Value = "Sarthak"
$ws.Cells.Item(4, 1).Value = "Sarthak2"
$ws.Cells.Item(5, 1).Value = "Arpit"
$ws.Cells.Item(6, 1).Value = "Harsh"
$ws.Cells.Item(7, 1).Value = "Jagdish"
$ws.Cells.Item(8, 1).Value = "Sushant"
$ws.Cells.Item(9, 1).Value = "Pranay"
$ws.Cells.Item(6, 1).Value = "Jagdish"
$ws.Cells.Item(7, 1).Value = "Harsh"

# --- riderID values in the new column B (rows 2-9)
$ws.Cells.Item(2, 2).Value = "rider1"
$ws.Cells.Item(3, 2).Value = "rider2"
$ws.Cells.Item(4, 2).Value = "rider3"
$ws.Cells.Item(5, 2).Value = "rider4"
$ws.Cells.Item(6, 2).Value = "rider5"
$ws.Cells.Item(7, 2).Value = "rider6"
$ws.Cells.Item(8, 2).Value = "rider7"
$ws.Cells.Item(9, 2).Value = "rider8"

# --- Header row
$ws.Cells.Item(1, 1).Value = "name"
$ws.Cells.Item(1, 3).Value = "phoneNumber"
$ws.Cells.Item(1, 4).Value = "totalBagVolume"
$ws.Cells.Item(1, 5).Value = "currentAvailableBagVolume"
$ws.Cells.Item(1, 2).Value = "riderID"

# --- Phone numbers, column C
$ws.Cells.Item(2, 3).Value = 9511725963
$ws.Cells.Item(3, 3).Value = 6239803560
$ws.Cells.Item(4, 3).Value = 1234567890
$ws.Cells.Item(5, 3).Value = 3456789123
$ws.Cells.Item(6, 3).Value = 3457899322
$ws.Cells.Item(7, 3).Value = 9680518959
$ws.Cells.Item(8, 3).Value = 6804188859
$ws.Cells.Item(9, 3).Value = 9680518923

# --- totalBagVolume / currentAvailableBagVolume formulas (columns D & E).
# Order matches the authored shared-formula group numbering: D2/E2 first,
# then the D3:D5 group, then the E3:E9 group, then D6/E6, then D7:D9 (E7:E9
# simply extends the existing E-group).
$ws.Range("D2").Formula = "=80*80*100"
$ws.Range("E2").Formula = "=D2"

$ws.Range("D3:D5").Formula = "=80*80*100"
$ws.Range("E3:E9").Formula = "=D3"

$ws.Range("D6").Formula = "=60*60*100"

$ws.Range("D7:D9").Formula = "=60*60*100"

# --- Column widths (D and E keep the widths inherited from the pre-insert
# C and D columns, so only A, B and C need to be resized explicitly)
$ws.Range("A1:B1").ColumnWidth = 31.5
$ws.Range("C1").ColumnWidth = 22.65

# --- Selection
$ws.Range("B2").Select()
